# Auto-generated script applying scheduled market-price refresh
# to the per-sheet Leve profit tables (columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 90.111115
$ws.Range("I6").Value = 101.25
$ws.Range("K6").Value = 303.75
$ws.Range("M6").Value = -191.75
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H88").Value = 3038033.8
$ws.Range("I88").Value = 11124274
$ws.Range("J88").Value = 5693.5
$ws.Range("K88").Value = 11124274
$ws.Range("L88").Value = 5693.5
$ws.Range("M88").Value = -11123868
$ws.Range("N88").Value = -6505.5
$ws.Range("H91").Value = 3038033.8
$ws.Range("I91").Value = 11124274
$ws.Range("J91").Value = 5693.5
$ws.Range("K91").Value = 11124274
$ws.Range("L91").Value = 5693.5
$ws.Range("M91").Value = -11122870
$ws.Range("N91").Value = -8501.5
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -64992
$ws.Range("H112").Value = 1927.7693
$ws.Range("I112").Value = 1122
$ws.Range("K112").Value = 3366
$ws.Range("M112").Value = -2258
$ws.Range("H113").Value = 25644080
$ws.Range("J113").Value = 3591
$ws.Range("L113").Value = 3591
$ws.Range("N113").Value = -10099

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2031.1428
$ws.Range("I32").Value = 2087.775
$ws.Range("K32").Value = 2087.775
$ws.Range("M32").Value = -1800.775
$ws.Range("H132").Value = 5357.086
$ws.Range("I132").Value = 4236.1055
$ws.Range("J132").Value = 6688.25
$ws.Range("K132").Value = 12708.3165
$ws.Range("L132").Value = 20064.75
$ws.Range("M132").Value = -10178.3165
$ws.Range("N132").Value = -25124.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14289070
$ws.Range("I20").Value = 17860486
$ws.Range("K20").Value = 17860486
$ws.Range("M20").Value = -17860239
$ws.Range("H86").Value = 2573
$ws.Range("I86").Value = 2128.4443
$ws.Range("K86").Value = 2128.4443
$ws.Range("M86").Value = -1005.4443
$ws.Range("H89").Value = 2573
$ws.Range("I89").Value = 2128.4443
$ws.Range("K89").Value = 10642.2215
$ws.Range("M89").Value = -5026.2215
$ws.Range("H105").Value = 9287640
$ws.Range("I105").Value = 626720.5
$ws.Range("J105").Value = 20835532
$ws.Range("K105").Value = 626720.5
$ws.Range("L105").Value = 20835532
$ws.Range("M105").Value = -624973.5
$ws.Range("N105").Value = -20839026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 458.9
$ws.Range("I7").Value = 273.375
$ws.Range("J7").Value = 1201
$ws.Range("K7").Value = 273.375
$ws.Range("L7").Value = 1201
$ws.Range("M7").Value = -160.375
$ws.Range("N7").Value = -1427
$ws.Range("H16").Value = 1964.2
$ws.Range("I16").Value = 1936.3334
$ws.Range("K16").Value = 1936.3334
$ws.Range("M16").Value = -1649.3334
$ws.Range("H31").Value = 5545.515
$ws.Range("I31").Value = 4735
$ws.Range("J31").Value = 6406.6875
$ws.Range("K31").Value = 4735
$ws.Range("L31").Value = 6406.6875
$ws.Range("M31").Value = -4440
$ws.Range("N31").Value = -6996.6875
$ws.Range("H34").Value = 5545.515
$ws.Range("I34").Value = 4735
$ws.Range("J34").Value = 6406.6875
$ws.Range("K34").Value = 4735
$ws.Range("L34").Value = 6406.6875
$ws.Range("M34").Value = -4533
$ws.Range("N34").Value = -6810.6875
$ws.Range("H58").Value = 3798.8333
$ws.Range("I58").Value = 3396.5
$ws.Range("K58").Value = 3396.5
$ws.Range("M58").Value = -3193.5
$ws.Range("H62").Value = 33338000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 33338000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H88").Value = 18214
$ws.Range("J88").Value = 18214
$ws.Range("L88").Value = 18214
$ws.Range("N88").Value = -19026
$ws.Range("H91").Value = 18214
$ws.Range("J91").Value = 18214
$ws.Range("L91").Value = 18214
$ws.Range("N91").Value = -21022
$ws.Range("H99").Value = 1789.8
$ws.Range("I99").Value = 1775
$ws.Range("J99").Value = 1799.6666
$ws.Range("K99").Value = 1775
$ws.Range("L99").Value = 1799.6666
$ws.Range("M99").Value = -277
$ws.Range("N99").Value = -4795.6666
$ws.Range("H113").Value = 1964.2
$ws.Range("I113").Value = 1936.3334
$ws.Range("K113").Value = 1936.3334
$ws.Range("M113").Value = 233.6666
$ws.Range("H120").Value = 32399.6
$ws.Range("J120").Value = 32399.6
$ws.Range("L120").Value = 32399.6
$ws.Range("N120").Value = -39657.6
$ws.Range("H126").Value = 1789.8
$ws.Range("I126").Value = 1775
$ws.Range("J126").Value = 1799.6666
$ws.Range("K126").Value = 5325
$ws.Range("L126").Value = 5398.9998
$ws.Range("M126").Value = -2855
$ws.Range("N126").Value = -10338.9998
$ws.Range("H132").Value = 9010769
$ws.Range("I132").Value = 1322.5714
$ws.Range("J132").Value = 37040160
$ws.Range("K132").Value = 3967.7142
$ws.Range("L132").Value = 111120480
$ws.Range("M132").Value = -1437.7142
$ws.Range("N132").Value = -111125540
$ws.Range("H136").Value = 3798.8333
$ws.Range("I136").Value = 3396.5
$ws.Range("K136").Value = 10189.5
$ws.Range("M136").Value = -7639.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31116044
$ws.Range("J4").Value = 14070202
$ws.Range("L4").Value = 42210606
$ws.Range("N4").Value = -42210830
$ws.Range("H13").Value = 347.8
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H44").Value = 2763.2727
$ws.Range("I44").Value = 279.4
$ws.Range("J44").Value = 4833.1665
$ws.Range("K44").Value = 838.1999999999999
$ws.Range("L44").Value = 14499.4995
$ws.Range("M44").Value = -440.1999999999999
$ws.Range("N44").Value = -15295.4995
$ws.Range("H132").Value = 2877.4
$ws.Range("J132").Value = 3665.6667
$ws.Range("L132").Value = 32991.0003
$ws.Range("N132").Value = -38051.0003
$ws.Range("H137").Value = 4928.143
$ws.Range("I137").Value = 4201.75
$ws.Range("J137").Value = 5218.7
$ws.Range("K137").Value = 12605.25
$ws.Range("L137").Value = 15656.1
$ws.Range("M137").Value = -7505.25
$ws.Range("N137").Value = -25856.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3779.8965
$ws.Range("J122").Value = 7636
$ws.Range("L122").Value = 22908
$ws.Range("N122").Value = -27808

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 749.8461
$ws.Range("I16").Value = 750.73914
$ws.Range("K16").Value = 750.73914
$ws.Range("M16").Value = -580.73914
$ws.Range("H54").Value = 19999
$ws.Range("I54").Value = 19999
$ws.Range("K54").Value = 19999
$ws.Range("M54").Value = -19355
$ws.Range("H64").Value = 38024.5
$ws.Range("J64").Value = 38024.5
$ws.Range("L64").Value = 38024.5
$ws.Range("N64").Value = -38474.5
$ws.Range("H67").Value = 38024.5
$ws.Range("J67").Value = 38024.5
$ws.Range("L67").Value = 38024.5
$ws.Range("N67").Value = -39584.5
$ws.Range("H94").Value = 165165
$ws.Range("J94").Value = 165165
$ws.Range("L94").Value = 165165
$ws.Range("N94").Value = -166517
$ws.Range("H132").Value = 7690.278
$ws.Range("I132").Value = 5043.8
$ws.Range("K132").Value = 15131.4
$ws.Range("M132").Value = -12601.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1020.4375
$ws.Range("I107").Value = 951
$ws.Range("K107").Value = 2853
$ws.Range("M107").Value = -933
$ws.Range("H113").Value = 1041.125
$ws.Range("I113").Value = 897.6923
$ws.Range("K113").Value = 2693.0769
$ws.Range("M113").Value = -523.0769
$ws.Range("H122").Value = 22729510
$ws.Range("I122").Value = 2698.8572
$ws.Range("J122").Value = 62501430
$ws.Range("K122").Value = 8096.571599999999
$ws.Range("L122").Value = 187504290
$ws.Range("M122").Value = -5646.571599999999
$ws.Range("N122").Value = -187509190
